$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 25246
$ws.Range("I16").Value = 25246
$ws.Range("K16").Value = 25246
$ws.Range("M16").Value = -25016

$ws.Range("H40").Value = 2047.4546
$ws.Range("I40").Value = 2166.3333
$ws.Range("J40").Value = 2002.875
$ws.Range("K40").Value = 2166.3333
$ws.Range("L40").Value = 2002.875
$ws.Range("M40").Value = -1991.3333
$ws.Range("N40").Value = -2352.875

$ws.Range("H64").Value = 5172
$ws.Range("I64").Value = 5172
$ws.Range("K64").Value = 5172
$ws.Range("M64").Value = -4924

$ws.Range("H67").Value = 5172
$ws.Range("I67").Value = 5172
$ws.Range("K67").Value = 5172
$ws.Range("M67").Value = -4314

$ws.Range("H69").Value = 18610.111
$ws.Range("I69").Value = 11499.5
$ws.Range("K69").Value = 34498.5
$ws.Range("M69").Value = -33624.5

$ws.Range("H72").Value = 18610.111
$ws.Range("I72").Value = 11499.5
$ws.Range("K72").Value = 103495.5
$ws.Range("M72").Value = -99127.5

$ws.Range("H86").Value = 50003052
$ws.Range("I86").Value = 58826676
$ws.Range("K86").Value = 58826676
$ws.Range("M86").Value = -58825553

$ws.Range("H88").Value = 3213.9285
$ws.Range("J88").Value = 4373
$ws.Range("L88").Value = 4373
$ws.Range("N88").Value = -5185

$ws.Range("H89").Value = 50003052
$ws.Range("I89").Value = 58826676
$ws.Range("K89").Value = 294133380
$ws.Range("M89").Value = -294127764

$ws.Range("H91").Value = 3213.9285
$ws.Range("J91").Value = 4373
$ws.Range("L91").Value = 4373
$ws.Range("N91").Value = -7181

$ws.Range("H136").Value = 147245
$ws.Range("J136").Value = 147245
$ws.Range("L136").Value = 147245
$ws.Range("N136").Value = -157445

$ws.Range("H138").Value = 5790.7676
$ws.Range("I138").Value = 12072.071
$ws.Range("K138").Value = 36216.213
$ws.Range("M138").Value = -31076.213

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 604077.5
$ws.Range("J74").Value = 930615.9
$ws.Range("L74").Value = 930615.9
$ws.Range("N74").Value = -932363.9

$ws.Range("H77").Value = 604077.5
$ws.Range("J77").Value = 930615.9
$ws.Range("L77").Value = 4653079.5
$ws.Range("N77").Value = -4661815.5

$ws.Range("H132").Value = 2265.9363
$ws.Range("I132").Value = 1979.875
$ws.Range("J132").Value = 3900.5715
$ws.Range("K132").Value = 5939.625
$ws.Range("L132").Value = 11701.7145
$ws.Range("M132").Value = -3409.625
$ws.Range("N132").Value = -16761.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 81247
$ws.Range("J135").Value = 81247
$ws.Range("L135").Value = 81247
$ws.Range("N135").Value = -91387

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2955862.2
$ws.Range("I6").Value = 3612476.2
$ws.Range("K6").Value = 3612476.2
$ws.Range("M6").Value = -3612363.2

$ws.Range("H16").Value = 15875270
$ws.Range("I16").Value = 28572526
$ws.Range("K16").Value = 28572526
$ws.Range("M16").Value = -28572239

$ws.Range("H31").Value = 3091.5952
$ws.Range("J31").Value = 2682.8484
$ws.Range("L31").Value = 2682.8484
$ws.Range("N31").Value = -3272.8484

$ws.Range("H34").Value = 3091.5952
$ws.Range("J34").Value = 2682.8484
$ws.Range("L34").Value = 2682.8484
$ws.Range("N34").Value = -3086.8484

$ws.Range("H86").Value = 31324.75
$ws.Range("I86").Value = 9005
$ws.Range("J86").Value = 38764.668
$ws.Range("K86").Value = 9005
$ws.Range("L86").Value = 38764.668
$ws.Range("M86").Value = -7882
$ws.Range("N86").Value = -41010.668

$ws.Range("H89").Value = 31324.75
$ws.Range("I89").Value = 9005
$ws.Range("J89").Value = 38764.668
$ws.Range("K89").Value = 45025
$ws.Range("L89").Value = 193823.34
$ws.Range("M89").Value = -39409
$ws.Range("N89").Value = -205055.34

$ws.Range("H113").Value = 15875270
$ws.Range("I113").Value = 28572526
$ws.Range("K113").Value = 28572526
$ws.Range("M113").Value = -28570356

$ws.Range("H122").Value = 2597.4443
$ws.Range("I122").Value = 2734.8096
$ws.Range("J122").Value = 2116.6667
$ws.Range("K122").Value = 8204.4288
$ws.Range("L122").Value = 6350.000100000001
$ws.Range("M122").Value = -5754.4288
$ws.Range("N122").Value = -11250.0001

$ws.Range("H132").Value = 45986.523
$ws.Range("I132").Value = 57799.168
$ws.Range("J132").Value = 3461
$ws.Range("K132").Value = 173397.504
$ws.Range("L132").Value = 10383
$ws.Range("M132").Value = -170867.504
$ws.Range("N132").Value = -15443

$ws.Range("H133").Value = 83119.664
$ws.Range("J133").Value = 83119.664
$ws.Range("L133").Value = 83119.664
$ws.Range("N133").Value = -88179.664

$ws.Range("H134").Value = 1478.2
$ws.Range("I134").Value = 1164.8462
$ws.Range("K134").Value = 3494.5386
$ws.Range("M134").Value = -959.5385999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2095
$ws.Range("I137").Value = 2095
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6285
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1185
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7641
$ws.Range("J46").Value = 6237.4
$ws.Range("L46").Value = 6237.4
$ws.Range("N46").Value = -6613.4

$ws.Range("H132").Value = 1555.6765
$ws.Range("J132").Value = 3071.889
$ws.Range("L132").Value = 9215.667000000001
$ws.Range("N132").Value = -14275.667

$ws.Range("H136").Value = 1268.6118
$ws.Range("I136").Value = 2050.0454
$ws.Range("J136").Value = 995.73016
$ws.Range("K136").Value = 6150.1362
$ws.Range("L136").Value = 2987.19048
$ws.Range("M136").Value = -3600.1362
$ws.Range("N136").Value = -8087.190479999999

$ws.Range("H139").Value = 150826.67
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 150826.67
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 150826.67
$ws.Range("N139").Value = -161106.67
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 18531.334
$ws.Range("J14").Value = 21510.334
$ws.Range("L14").Value = 21510.334
$ws.Range("N14").Value = -21846.334

$ws.Range("H122").Value = 2458.238
$ws.Range("I122").Value = 1708.4706
$ws.Range("K122").Value = 5125.4118
$ws.Range("M122").Value = -2675.4118

$ws.Range("H126").Value = 2606.3125
$ws.Range("I126").Value = 2323.3635
$ws.Range("K126").Value = 6970.0905
$ws.Range("M126").Value = -4500.0905

$ws.Range("H132").Value = 1621.3556
$ws.Range("I132").Value = 991.7826
$ws.Range("K132").Value = 2975.3478
$ws.Range("M132").Value = -445.3478
